# Add descriptions/titles for the DMI Reference Distributeur extension.
#
# - Metadata sheet: fill in the Title and Description values for this
#   StructureDefinition, and bump the generation Date.
# - Elements sheet: give the root "Extension" row a proper Short / Definition
#   (replacing the generic "Optional Extensions Element" boilerplate) and
#   clear its now-irrelevant RIM mapping.

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B5").Value = "DMI Reference Distributeur"
$wsMeta.Range("B8").Value = "2026-02-25T08:15:31+00:00"
$wsMeta.Range("B12").Value = "Extension créée dans ce volet pour représenter la référence distributeur."

$wsElem = $wb.Worksheets.Item("Elements")
$wsElem.Range("L2").Value = "DMI Reference Distributeur"
$wsElem.Range("M2").Value = "Extension créée dans ce volet pour représenter la référence distributeur."
$wsElem.Range("AK2").Value = ""
